# Commit: Change '_' to '-' in Excel metadata (Close #5)
#
# The PX-file metadata keywords stored on the "General_MD" sheet used an
# underscore between the keyword "stem" and its language suffix
# (e.g. SUBJECT_AREA_en). They should use a hyphen before the language
# suffix instead (SUBJECT-AREA_en), matching the PX-file keyword spec.
#
# NB: the order below intentionally matches the order the renamed keywords
# were (re)written in, since Excel appends newly-typed strings to the
# shared-string table in write order.

$wb = $excel.ActiveWorkbook
$general = $wb.Worksheets.Item("General_MD")

$general.Range("A3").Value  = "AXIS-VERSION"
$general.Range("A19").Value = "SUBJECT-CODE"
$general.Range("A24").Value = "CREATION-DATE"
$general.Range("A25").Value = "UPDATE-FREQUENCY"
$general.Range("A26").Value = "LAST-UPDATED"
$general.Range("A27").Value = "NEXT-UPDATE"
$general.Range("A20").Value = "SUBJECT-AREA_fo"
$general.Range("A21").Value = "SUBJECT-AREA_en"

# --- Sheet selections / active tab -----------------------------------
# Variables_MD: move the lingering selection to K1 (no tab change).
$variables = $wb.Worksheets.Item("Variables_MD")
[void]$variables.Range("K1").Select()

# Codelists_2MD: selection stays on C12; it just stops being the active tab.
$codelists = $wb.Worksheets.Item("Codelists_2MD")
[void]$codelists.Range("C12").Select()

# General_MD becomes the active/selected tab, with A7 selected.
[void]$general.Range("A7").Select()
